$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.736.10'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.421.04'
$ws.Range("E3").Value = '  -2.33%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.84'
$ws.Range("E5").Value = '  -1.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.73'
$ws.Range("E6").Value = '  -3.43%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.480'
$ws.Range("E8").Value = '  -1.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.57'
$ws.Range("E9").Value = '  +3.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("E10").Value = '  +0.85%  '

$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.016.56'
$ws.Range("E12").Value = '  -2.08%  '

$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000177'
$ws.Range("E14").Value = '  -1.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.429.41'
$ws.Range("E15").Value = '  -2.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.736.01'
$ws.Range("E16").Value = '  -0.91%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.28'
$ws.Range("E17").Value = '  -1.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.85'
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.67'
$ws.Range("E19").Value = '  -1.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.36'
$ws.Range("E20").Value = '  -1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '384.02'
$ws.Range("E21").Value = '  -2.29%  '

$ws.Range("E22").Value = '  -1.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.563.26'
$ws.Range("E23").Value = '  -2.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.14'
$ws.Range("E24").Value = '  -0.61%  '

$ws.Range("E25").Value = '  -0.26%  '

$ws.Range("E26").Value = '  -4.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  -1.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.03'
$ws.Range("E29").Value = '  -4.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.94'
$ws.Range("E30").Value = '  -3.49%  '

$ws.Range("E31").Value = '  +1.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.41'
$ws.Range("E32").Value = '  -3.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.452.99'
$ws.Range("E33").Value = '  -2.11%  '

$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.84'
$ws.Range("E35").Value = '  -2.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.16'
$ws.Range("E36").Value = '  +0.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.74'
$ws.Range("E37").Value = '  -1.83%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.52'
$ws.Range("E38").Value = '  -1.64%  '

$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.00'
$ws.Range("E39").Value = '  -2.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0772'
$ws.Range("E40").Value = '  -0.85%  '

$ws.Range("E41").Value = '  -2.69%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.43'
$ws.Range("E43").Value = '  -1.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.34'
$ws.Range("E44").Value = '  -1.28%  '

$ws.Range("E45").Value = '  -2.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.30'
$ws.Range("E46").Value = '  -8.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.10'
$ws.Range("E47").Value = '  -5.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.72'
$ws.Range("E48").Value = '  -0.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.890'
$ws.Range("E49").Value = '  +0.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.282.19'
$ws.Range("E50").Value = '  -2.54%  '

$ws.Range("E51").Value = '  -2.69%  '
